# Update the "About" sheet: two explanatory sentences swap places (A8 <-> A9)
$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$about.Range("A8").Value = "for the United States.  We arbitrarily assign priority 2 to all other plant types."
$about.Range("A9").Value = "only types for which a non-zero quantity is specified for guaranteed dispatch in the BAU case"

# Update the "BDPbES" sheet: electricity source subscript gains/loses members
$ws = $wb.Worksheets.Item("BDPbES")

# Row 3 "natural gas steam turbine" is renamed to "natural gas nonpeaker"
$ws.Range("A3").Value = "natural gas nonpeaker"

# Row 4 "natural gas combined cycle" is removed outright (rows below shift up one)
$ws.Rows.Item(4).Delete()

# Five new electricity sources are appended at the bottom of the table (now rows 18-23)
$newSources = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$row = 18
foreach ($name in $newSources) {
    $ws.Range("A$row").Value = $name
    $ws.Range("B$row").Value = 2
    $ws.Range("C$row`:AK$row").Formula = '=$B' + $row
    $row = $row + 1
}

$ws.Columns.Item(1).ColumnWidth = 22.7109375
$ws.Range("A24").Select()
